# s4_cover.docx — "updated scrum master for s4 cover"
#
# 1. The first three (empty) title-block paragraphs have their sole,
#    unformatted run promoted to carry the same direct character
#    formatting (Times New Roman / 36 half-pt-squared size / single
#    underline) that the paragraph mark already carries.
# 2. The title paragraph's three runs ("The Scrumbags – Sprint" / "4" /
#    " Deliverable") collapse into a single run with the full text.
# 3. Nathan Ackerman becomes the Scrum Master instead of Wei-Hao Chen:
#    "Nathan Ackerman – NathanAckerman" gains a "(Scrum Master) " run in
#    the middle (split across 3 runs), and Wei-Hao Chen's line drops the
#    "(Scrum Master)" suffix.
# 4. The date line's two runs ("19" / ", 2017") collapse into one.

$d = $word.ActiveDocument
$EnDash = [char]0x2013

# --- 1. format the three empty paragraphs at the top of the cover ---
for ($i = 1; $i -le 3; $i++) {
    $r = $d.Paragraphs($i).Range
    $r.Font.Name   = "Times New Roman"
    $r.Font.NameBi = "Times New Roman"
    $r.Font.Size   = 18
    $r.Font.SizeBi = 18
    $r.Font.Underline = 1
}

# --- 2. merge the title runs into a single run ---
$title = "The Scrumbags " + $EnDash + " Sprint 4 Deliverable"
$d.Content.Find.Execute($title, $true, $false, $false, $false, $false, $true, 1, $false, $title, 2)

# --- 3a. Wei-Hao Chen line: drop the "(Scrum Master)" title ---
$weiOld = "Wei-Hao Chen (Scrum Master)  " + $EnDash + " yoshino0705"
$weiNew = "Wei-Hao Chen " + $EnDash + " yoshino0705"
$d.Content.Find.Execute($weiOld, $true, $false, $false, $false, $false, $true, 1, $false, $weiNew, 2)

# --- 3b. Nathan Ackerman line: add "(Scrum Master) " in the middle ---
$nathanOld = "Nathan Ackerman " + $EnDash + " NathanAckerman"
$nathanNew = "Nathan Ackerman  (Scrum Master) " + $EnDash + " NathanAckerman"
$d.Content.Find.Execute($nathanOld, $true, $false, $false, $false, $false, $true, 1, $false, $nathanNew, 2)

# Locate the (now rewritten) Nathan Ackerman paragraph and split it into
# three runs: "Nathan Ackerman  " / "(Scrum Master) " / "– NathanAckerman"
$nathanPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Nathan Ackerman")) {
        $nathanPara = $p
    }
}
$pStart = $nathanPara.Range.Start
$cut1 = $pStart + 17   # end of "Nathan Ackerman  "
$cut2 = $pStart + 32   # end of "(Scrum Master) "

$runA = $d.Range($pStart, $cut1)
$runB = $d.Range($cut1, $cut2)

# Toggling a direct-formatting property on and back forces the run
# boundary to persist even though the resulting formatting is identical
# to its neighbours.
$runA.Bold = 1
$runA.Bold = 0
$runB.Bold = 1
$runB.Bold = 0

# --- 4. merge the date runs ("19" + ", 2017") into one ---
$d.Content.Find.Execute("19, 2017", $true, $false, $false, $false, $false, $true, 1, $false, "19, 2017", 2)
